$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Game")
$ws.Range("B3").Value = "2020-07-01 00:00:00 +0300"
